$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps text formatting so numeric-looking strings
# are not reinterpreted as numbers.
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range("D2").Value = "247.84"
$ws.Range("D4").Value = "5.509"
$ws.Range("D5").Value = "0.05623"
$ws.Range("D7").Value = "0.8039"
$ws.Range("D8").Value = "1.050"
$ws.Range("D9").Value = "0.1424"
$ws.Range("D10").Value = "0.07312"
$ws.Range("D11").Value = "0.03191"
$ws.Range("D12").Value = "0.02937"
$ws.Range("D13").Value = "0.09270"
$ws.Range("D14").Value = "0.001673"
$ws.Range("D15").Value = "3.208"
$ws.Range("D16").Value = "0.04690"
$ws.Range("D17").Value = "0.0005918"
$ws.Range("E17").Value = "16OneONEWorstin24h"
$ws.Range("D18").Value = "0.006273"
$ws.Range("D19").Value = "0.001056"
$ws.Range("D20").Value = "0.003822"
$ws.Range("D22").Value = "0.0003305"
$ws.Range("D24").Value = "3.384"
$ws.Range("D25").Value = "2.093"
$ws.Range("D40").Value = "0.04162"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006870"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "0.003505"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "0.1042"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "0.009054"
$ws.Range("D45").Value = "0.00005640"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").Value = "0.6810"
$ws.Range("D48").Value = "0.02525"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("D50").Value = "0.01011"

# Reset style attribute on Price column so it matches the original (default) style,
# since setting NumberFormat alone would otherwise leave a style index on the cells.
$ws.Range("D2:D50").Style = "Normal"
